$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update two existing result cells ---
# Tabela3 (Modularity results): Girvan-Newman column, Email row
$ws.Range("V3").Value = 0.422
# Tabela4 (Accuracy results): Infomap column, Email row
$ws.Range("V10").Value = 0.6161

# --- Add the new "Times (s)" mini table (Tabela46) header row ---
$ws.Range("R12").Value = "Times (s)"
$ws.Range("S12").Value = "Louvain"
$ws.Range("T12").Value = "Leiden"
$ws.Range("U12").Value = "Girvan-Newman"
$ws.Range("V12").Value = "Infomap"

# --- Add the new "Times (s)" mini table (Tabela46) data row ---
$ws.Range("R13").Value = "Email"
$ws.Range("S13").Value = 0.331
$ws.Range("T13").Value = 0.07
$ws.Range("U13").Value = 382.82
$ws.Range("V13").Value = 0.261

# --- Turn R12:V13 into a proper Excel Table ("Tabela46") ---
$lo = $ws.ListObjects.Add(1, $ws.Range("R12:V13"), 0, 1)
$lo.Name = "Tabela46"
$lo.TableStyle = "TableStyleLight13"

# Reuse the exact same cell formatting already used by the sibling
# "Tabela4" header/data rows (row 9/10) and "Tabela3" data row (row 3)
# so the new cells share the same look (bold first column, thin borders,
# centered text, "0.000" number format for the numeric results).
$ws.Range("R9").Copy()
$ws.Range("R12").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("S9").Copy()
$ws.Range("S12:U12").PasteSpecial(-4122)
$ws.Range("V9").Copy()
$ws.Range("V12").PasteSpecial(-4122)

$ws.Range("R3").Copy()
$ws.Range("R13").PasteSpecial(-4122)
$ws.Range("S3").Copy()
$ws.Range("S13:U13").PasteSpecial(-4122)
$ws.Range("V3").Copy()
$ws.Range("V13").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# --- Restore the selection saved in the authored workbook ---
$ws.Range("U14").Select()
